$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null

$ws.Range("G3").Copy() | Out-Null
$ws.Range("G8").PasteSpecial(-4122) | Out-Null

$ws.Range("A8").Value = 42649.656481481485
$ws.Range("B8").Value = $true
$ws.Range("C8").Value = 10043.799999999999
$ws.Range("D8").Value = 9993.33
$ws.Range("E8").Value = 18.829999999999998
$ws.Range("F8").Value = 19.02
$ws.Range("G8").Value = $false
$ws.Range("H8").Value = 1.01
$ws.Range("I8").Value = $false
